$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.370.03'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.605.14'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0854'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '1.830.70'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').Value = '1.605.34'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.28'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '26.374.19'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.87'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.58%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  +4.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('E23').Value = '  +2.80%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.39'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = '1.493.04'
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.562'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.819'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.929'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.14%  '
$ws.Range('D44').Value = '1.743.72'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.759'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.83'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.98%  '
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0959'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +0.02%  '
